$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("T6").Value  = [double]"426.37093227978386"
$ws.Range("T7").Value  = [double]"1.2134608580999326"
$ws.Range("T8").Value  = [double]"182.09200851950402"
$ws.Range("T9").Value  = [double]"-5.9529071385460965E-2"
$ws.Range("T10").Value = [double]"0.96150950769486943"
$ws.Range("T11").Value = [double]"1.3944631257729466E-4"
$ws.Range("T12").Value = [double]"-4.1526492641042457E-5"
$ws.Range("T13").Value = [double]"-8.0668440966409403E-8"
$ws.Range("T14").Value = [double]"0.94977136273332075"
$ws.Range("T15").Value = [double]"1.2004675011055647E-4"
$ws.Range("T16").Value = [double]"-0.15528480952384874"
$ws.Range("T17").Value = [double]"-7.1647109472057303E-8"
$ws.Range("T18").Value = [double]"1.2021082917533591E-4"
$ws.Range("T19").Value = [double]"-0.41490204204113817"

$excel.ActiveWindow.Left = 3394
$excel.ActiveWindow.Top = 3394
